$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark (currently sitting right after
#    the "Navigation" paragraph's text).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Remove the whole "MatFile Reader einbinden (...)" list paragraph (the
#    last paragraph in the body, right after "Unit Test").
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.Delete()

# 3. Re-create the "_GoBack" bookmark at the end of the "Unit Test"
#    paragraph's text (now the last paragraph), matching the post-edit
#    document. A zero-length Range placed exactly on a paragraph's last
#    character position trips up this host's Bookmarks.Add, so work around
#    it: temporarily append a placeholder character after "Unit Test",
#    plant the (now safely non-boundary) bookmark in front of it, then
#    remove the placeholder again - the bookmark stays put, collapsed,
#    right before the paragraph mark.
$unitTest = $d.Paragraphs($d.Paragraphs.Count)
$insPos = $unitTest.Range.End - 1
$placeholder = $d.Range($insPos, $insPos)
$placeholder.InsertAfter("Z")

$bmPos = $insPos
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholderRange = $d.Range($bmPos, $bmPos + 1)
$placeholderRange.Delete()
